$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update July 20 data for Los Angeles Lakers (C), Utah Jazz (I), New Orleans Pelicans (S)
# and recolor the refreshed cells with the new font color (#333333).
$ws.Range("C2").Value2 = 100
$ws.Range("I2").Value2 = 103.9
$ws.Range("S2").Value2 = 103
$ws.Range("C3").Value2 = 99
$ws.Range("I3").Value2 = 103
$ws.Range("S3").Value2 = 103.9
$ws.Range("C4").Value2 = 1
$ws.Range("I4").Value2 = 1
$ws.Range("S4").Value2 = -1
$ws.Range("C5").Value2 = 65.599999999999994
$ws.Range("I5").Value2 = 45.9
$ws.Range("S5").Value2 = 55
$ws.Range("C6").Value2 = 1.31
$ws.Range("I6").Value2 = 0.85
$ws.Range("S6").Value2 = 1.05
$ws.Range("C7").Value2 = 15.2
$ws.Range("I7").Value2 = 12.5
$ws.Range("S7").Value2 = 15.3
$ws.Range("C8").Value2 = 25.5
$ws.Range("I8").Value2 = 34
$ws.Range("S8").Value2 = 36.799999999999997
$ws.Range("C9").Value2 = 81.400000000000006
$ws.Range("I9").Value2 = 63.2
$ws.Range("S9").Value2 = 66
$ws.Range("C10").Value2 = 50
$ws.Range("I10").Value2 = 49.5
$ws.Range("S10").Value2 = 50.5
$ws.Range("C11").Value2 = 15.5
$ws.Range("I11").Value2 = 19.600000000000001
$ws.Range("S11").Value2 = 20.8
$ws.Range("C12").Value2 = 45.7
$ws.Range("I12").Value2 = 48.8
$ws.Range("S12").Value2 = 48.4
$ws.Range("C13").Value2 = 52.4
$ws.Range("I13").Value2 = 55
$ws.Range("S13").Value2 = 51
$ws.Range("C14").Value2 = 102.5
$ws.Range("I14").Value2 = 101.5
$ws.Range("S14").Value2 = 101.5
$ws.Range("C15").Value2 = 53.7
$ws.Range("I15").Value2 = 52.8
$ws.Range("S15").Value2 = 47.2

$updatedRanges = @("C2:C15", "I2:I15", "S2:S15")
foreach ($rangeAddress in $updatedRanges) {
    $ws.Range($rangeAddress).Font.Color = 3355443
}

[void]$ws.Range("A17").Select()
